$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift all data one column to the right (A:D -> B:E) by inserting a new
# column A and pushing the existing content over.
$ws.Columns("A:A").Insert()

# Match the saved selection state: whole-column A selected, no active cell override.
$ws.Range("A1:A1048576").Select()
